$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Propagate the CURRENT ("yellow", fillId4) row format from row 9 down to
#    the three brand-new rows (13,14,15) BEFORE we repaint rows 9-12 white.
#    Doing this first lets the new rows inherit the exact pre-existing
#    "yellow" cell styles instead of creating fresh duplicate ones.
# ---------------------------------------------------------------------------
$ws.Range("A9:C9").Copy()
$ws.Range("A13:C15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Fill in the new row contents (row-major order so new shared strings are
#    appended in the same sequence as the target file: Logear usuario sin
#    JWT, Obtener rutas disponibles, the new reportes URL, then the excel
#    download description).
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "GET"
$ws.Range("B13").Value = "http://localhost:9090/usuarios"
$ws.Range("C13").Value = "Logear usuario sin JWT"

$ws.Range("A14").Value = "GET"
$ws.Range("B14").Value = "http://localhost:9090/rutas"
$ws.Range("C14").Value = "Obtener rutas disponibles"

$ws.Range("A15").Value = "GET"
$ws.Range("B15").Value = "http://localhost:9090/rutas/reportes"
$ws.Range("C15").Value = "Servicio que descarga excel con reporte total programados"

# ---------------------------------------------------------------------------
# 3) Hyperlink the three new endpoint cells. Hyperlinks.Add() stamps its own
#    (duplicate) style on the target cell, so immediately re-apply the
#    "yellow hyperlink" format that is already sitting on B9 to put the
#    style index back where it belongs.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B13"), "http://localhost:9090/usuarios") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B14"), "http://localhost:9090/rutas") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B15"), "http://localhost:9090/rutas/reportes") | Out-Null

$ws.Range("B9").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Rows 9-12 switch from the yellow highlight to a plain white (theme)
#    fill. Restyle row 9 directly, then fan that exact format out to rows
#    10-12 via copy/paste-formats so every row reuses the same new style
#    indices instead of minting one per row.
# ---------------------------------------------------------------------------
$ws.Range("A9").Interior.ThemeColor = 2
$ws.Range("B9").Interior.ThemeColor = 2
$ws.Range("C9").Interior.ThemeColor = 2

$ws.Range("A9:C9").Copy()
$ws.Range("A10:C12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5) Column C widens slightly to comfortably fit the new, longer
#    descriptions (target OOXML width ~53.57; 52.65 is the closest input
#    this engine's char<->pixel rounding can reach).
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 52.65

# ---------------------------------------------------------------------------
# 6) Match the final selection left behind by the edit.
# ---------------------------------------------------------------------------
$ws.Range("B23").Select() | Out-Null
